$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ticker cells that changed between before/after states
$ws.Range("B2").Value = "NSE:A2ZINFRA"
$ws.Range("C2").Value = "NSE:3IINFOLTD"
$ws.Range("D2").Value = "NSE:NYKAA"
$ws.Range("E2").Value = "NSE:ANGELONE"
$ws.Range("F2").Value = "NSE:CIPLA"
$ws.Range("B3").Value = "NSE:ANANTRAJ"
$ws.Range("C3").Value = "NSE:AARVI"
$ws.Range("E3").Value = "NSE:BRITANNIA"
$ws.Range("F3").Value = "NSE:NCC"
$ws.Range("B4").Value = "NSE:BANSWRAS"
$ws.Range("C4").Value = "NSE:ASHAPURMIN"
$ws.Range("E4").Value = "NSE:NTPC"
$ws.Range("F4").Value = "NSE:PAYTM"
$ws.Range("B5").Value = "NSE:BBETF0432"
$ws.Range("C5").Value = "NSE:AXISGOLD"
$ws.Range("E5").Value = "NSE:ONGC"
$ws.Range("F5").Value = "NSE:PRESTIGE"
$ws.Range("B6").Value = "NSE:BLUEJET"
$ws.Range("C6").Value = "NSE:CEATLTD"
$ws.Range("B7").Value = "NSE:CIPLA"
$ws.Range("C7").Value = "NSE:CENTURYPLY"
$ws.Range("B8").Value = "NSE:DOLLAR"
$ws.Range("C8").Value = "NSE:CYIENT"
$ws.Range("B9").Value = "NSE:DONEAR"
$ws.Range("C9").Value = "NSE:FORTIS"
$ws.Range("B10").Value = "NSE:HCG"
$ws.Range("C10").Value = "NSE:GREENPLY"
$ws.Range("B11").Value = "NSE:IDFNIFTYET"
$ws.Range("C11").Value = "NSE:HARDWYN"
$ws.Range("B12").Value = "NSE:INDIANCARD"
$ws.Range("C12").Value = "NSE:HGS"
$ws.Range("B13").Value = "NSE:JAGSNPHARM"
$ws.Range("C13").Value = "NSE:HNGSNGBEES"
$ws.Range("B14").Value = "NSE:MBLINFRA"
$ws.Range("C14").Value = "NSE:ICEMAKE"
$ws.Range("B15").Value = "NSE:MON100"
$ws.Range("C15").Value = "NSE:IRFC"
$ws.Range("B16").Value = "NSE:NAHARPOLY"
$ws.Range("C16").Value = "NSE:KALAMANDIR"
$ws.Range("B17").Value = "NSE:NAHARSPING"
$ws.Range("C17").Value = "NSE:MTNL"
$ws.Range("B18").Value = "NSE:PAYTM"
$ws.Range("C18").Value = "NSE:NIFTYQLITY"
$ws.Range("B19").Value = "NSE:PRICOLLTD"
$ws.Range("C19").Value = "NSE:OMAXAUTO"
$ws.Range("B20").Value = "NSE:SALZERELEC"
$ws.Range("C20").Value = "NSE:OSWALGREEN"
$ws.Range("C21").Value = "NSE:PGHL"
$ws.Range("C22").Value = "NSE:PRINCEPIPE"
$ws.Range("C23").Value = "NSE:RUSHIL"

# Remove the now-unused trailing rows (24-29) so the sheet dimension shrinks to A1:F23
$ws.Range("A24:F29").Delete()
